$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new value in cell B10
$ws.Range("B10").Value = "ebe"

# Move the active selection to N17, as it was left after the edit
$ws.Range("N17").Select()
